# Reorders the compound data rows (2-9) on the active sheet.
# Row 1 (headers) and row 5 (palmitic acid) are unchanged; the other
# rows are permuted so that each row ends up holding a different
# compound's data, per the two cycles:
#   row2 <- row7, row7 <- row9, row9 <- row4, row4 <- row2
#   row3 <- row6, row6 <- row8, row8 <- row3
#
# Implementation: snapshot every source row's values first (so later
# writes don't clobber data still needed), then write the snapshots
# into their destination rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 2
$lastRow = 9
$lastCol = 24   # column X

# New row (key) <- Old row (value) it should receive data from.
$mapping = @{
    2 = 7
    3 = 6
    4 = 2
    5 = 5
    6 = 8
    7 = 9
    8 = 3
    9 = 4
}

# Snapshot all source rows' values before writing anything.
$snapshots = @{}
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $rowValues = @()
    for ($c = 1; $c -le $lastCol; $c++) {
        $rowValues += $ws.Cells.Item($r, $c).Value()
    }
    $snapshots[$r] = $rowValues
}

# Write each destination row from its mapped source snapshot.
foreach ($destRow in $mapping.Keys) {
    $srcRow = $mapping[$destRow]
    $values = $snapshots[$srcRow]
    for ($c = 1; $c -le $lastCol; $c++) {
        $ws.Cells.Item($destRow, $c).Value = $values[$c - 1]
    }
}

$wb.Save()
